# Apply "updated combined results, added instances that dCAQE determines wrong result"
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("runs")
$ws2 = $wb.Worksheets.Item("realtime")

# --- "runs" sheet: rename the dCAQE run-set label (N2:Q2) -----------------
$ws1.Range("N2:Q2").Value = "qbf-mode-noproofreconstruction.SAT"

# --- "runs" sheet: updated dCAQE measurements (columns N-Q, rows 4-25) ----
$ws1.Range("O4").Value = 38.873337665000001
$ws1.Range("P4").Value = 38.877165973186401

$ws1.Range("O5").Value = 210.39101399399999
$ws1.Range("P5").Value = 210.26516573130999
$ws1.Range("Q5").Value = 36.704256000000001

$ws1.Range("O6").Value = 198.09790560499999
$ws1.Range("P6").Value = 198.008979242295
$ws1.Range("Q6").Value = 112.332799999999

$ws1.Range("N7").Value = "unsat"
$ws1.Range("O7").Value = 341.41489188700001
$ws1.Range("P7").Value = 341.19318745285199
$ws1.Range("Q7").Value = 39.297024

$ws1.Range("O8").Value = 901.69190069199999
$ws1.Range("P8").Value = 901.06110102310697

$ws1.Range("N9").Value = "unsat"
$ws1.Range("O9").Value = 364.90279952600002
$ws1.Range("P9").Value = 364.66176222264698
$ws1.Range("Q9").Value = 22.990848

$ws1.Range("O10").Value = 901.59020197200005
$ws1.Range("P10").Value = 901.04938426241199
$ws1.Range("Q10").Value = 32.014336

$ws1.Range("O11").Value = 901.69378331799999
$ws1.Range("P11").Value = 901.06163898110299
$ws1.Range("Q11").Value = 64.020479999999907

$ws1.Range("O12").Value = 901.69429098000001
$ws1.Range("P12").Value = 901.06575049832395
$ws1.Range("Q12").Value = 42.352640000000001

$ws1.Range("N13").Value = "unsat"
$ws1.Range("O13").Value = 290.70346587900002
$ws1.Range("P13").Value = 290.53796089067998
$ws1.Range("Q13").Value = 47.423487999999999

$ws1.Range("O14").Value = 901.69147759500004
$ws1.Range("P14").Value = 901.06580974906603

$ws1.Range("O15").Value = 570.91204480299996
$ws1.Range("P15").Value = 570.50507144629898

$ws1.Range("O16").Value = 901.69746372500003
$ws1.Range("P16").Value = 901.08165337145294

$ws1.Range("O17").Value = 901.593336357
$ws1.Range("P17").Value = 901.04965047910798
$ws1.Range("Q17").Value = 93.888511999999906

$ws1.Range("O18").Value = 901.69411709600001
$ws1.Range("P18").Value = 901.08165212720598
$ws1.Range("Q18").Value = 30.441471999999901

$ws1.Range("O19").Value = 901.69445660400004
$ws1.Range("P19").Value = 901.04966185241904
$ws1.Range("Q19").Value = 32.522239999999996

$ws1.Range("O20").Value = 901.69414579700003
$ws1.Range("P20").Value = 901.08171526342596

$ws1.Range("O21").Value = 901.68949501400004
$ws1.Range("P21").Value = 901.04978442564595

$ws1.Range("O22").Value = 901.69221876899996
$ws1.Range("P22").Value = 901.03381768241502

$ws1.Range("O23").Value = 901.57996135200005
$ws1.Range("P23").Value = 901.038075443357
$ws1.Range("Q23").Value = 43.081727999999998

$ws1.Range("N24").Value = "unsat"
$ws1.Range("O24").Value = 582.11026493600002
$ws1.Range("P24").Value = 581.70183198898997
$ws1.Range("Q24").Value = 27.328512

$ws1.Range("O25").Value = 901.67930187900004
$ws1.Range("P25").Value = 901.06177636235896
$ws1.Range("Q25").Value = 59.903999999999897

# --- "realtime" sheet: F4:F25 now pull the status live from runs!N -------
for ($r = 4; $r -le 25; $r++) {
    $ws2.Range("F$r").Formula = "=runs!N$r"
}

# --- View state: selection / scroll position ------------------------------
$ws1.Range("I1").Select()
$excel.ActiveWindow.ScrollColumn = 9
$ws1.Range("P29").Select()

$ws2.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$ws2.Range("H31").Select()
